$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# The only actual content change is cell E8: "Good Morning" -> "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the selection change recorded in the sheet view (active cell E8)
$ws.Range("E8").Select()
